# Insert a new "id" column at the very front of the sheet (everything
# that used to live in columns A:J shifts right to B:K), then fill the
# new column A with a header label and a 1-based running row index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing columns (A:J) one place to the right.
$ws.Columns("A:A").Insert()

# Header for the new id column.
$ws.Range("A1").Value = "id"

# Sequential id values for each data row (row 2 -> 1, row 3 -> 2, ...).
$lastRow = 35
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# The row that held the 1936 entry grew a touch taller after the edit.
$ws.Rows(12).RowHeight = 168.75

# Restore the view roughly where the author left it - scrolled further
# down the table with the selection near the bottom.
$ws.Range("C39").Select()
